$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume (E) updates per latest crypto snapshot.
# D-column values are forced to text (NumberFormat "@") so that
# numeric-looking strings like "598.23" are not reinterpreted as
# numbers by Excel; the cell style is then reset back to Normal
# so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.208.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.626.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.91%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  +2.94%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.625.16'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("E10").Value = '  +3.17%  '

$ws.Range("E11").Value = '  +0.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.348'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.103.99'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("E16").Value = '  +0.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.101.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.626.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '363.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.07%  '

$ws.Range("E21").Value = '  -3.45%  '

$ws.Range("E22").Value = '  -0.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.45%  '

$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '67.41'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.758.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.02'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000101'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '570.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.26%  '

$ws.Range("E31").Value = '  -3.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.59%  '

$ws.Range("E33").Value = '  -0.33%  '

$ws.Range("E34").Value = '  +0.05%  '

$ws.Range("E35").Value = '  -3.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.52'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.367'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("E41").Value = '  -3.16%  '

$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.18%  '

$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("E46").Value = '  -0.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '155.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.54%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0283'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.620'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.48%  '
